$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grab the existing teal/Consolas style (already used on A2/B2) before we
# clear the sheet, so we can re-apply the identical style later without
# Excel minting a brand-new (but visually-equivalent) font/style entry.
$styleSource = $ws.Range("A2")
$styleSource.Copy()
$ws.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats, parked far away from the real data
$excel.CutCopyMode = $false

# Clear out the old content area (A1:E10) before laying out the new table
$ws.Range("A1:E10").Clear()

# New header row
$ws.Cells.Item(1, 1).Value = "Feature Step Name"
$ws.Cells.Item(1, 2).Value = "Action"
$ws.Cells.Item(1, 3).Value = "Xpath"
$ws.Cells.Item(1, 4).Value = "Value"

# New data rows (feature step name, action, xpath, value)
$data = @(
    @("fill_all_mandatory_fields", "Sendkeys", "//*[@id=\`"firstname\`"]", "Automation 1"),
    @("fill_all_mandatory_fields", "Sendkeys", '//*[@id="lastname"]', "test"),
    @("fill_all_mandatory_fields", "Sendkeys", '//*[@id="email_address"]', "emailmaster@mailinator.com"),
    @("fill_all_mandatory_fields", "Sendkeys", '//*[@id="password"]', "Test@123"),
    @("fill_all_mandatory_fields", "Sendkeys", '//*[@id="password-confirmation"]', "Test@123"),
    @("fill_all_mandatory_fields", "Select", "//*[@id=\`"gender\`"]", "Male")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply the "teal Consolas" style (style index 1) to columns A, C and D of the
# data rows (rows 2-7) - matches xf with fontId=1 used in the original workbook.
foreach ($colRange in @("A2:A7", "C2:C7", "D2:D7")) {
    $ws.Range("Z1").Copy()
    $ws.Range($colRange).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Column widths (values chosen so the pixel-snapped COM width lands on the
# closest achievable match to the target stored widths of 30.140625 /
# 67.140625 / 47.42578125 / 35.42578125)
$ws.Columns.Item(1).ColumnWidth = 29.3
$ws.Columns.Item(2).ColumnWidth = 66.3
$ws.Columns.Item(3).ColumnWidth = 46.65
$ws.Columns.Item(4).ColumnWidth = 34.65

# Selection matches the saved state in the new workbook
$ws.Range("B6").Select()

# Update Sheet2 (same wording, listed here for completeness / safety, values unchanged)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Cells.Item(2, 1).Value = "operation"
$ws2.Cells.Item(3, 1).Value = "click"
$ws2.Cells.Item(3, 4).Value = "xpath"
$ws2.Cells.Item(4, 1).Value = "sendkeys"
$ws2.Cells.Item(5, 1).Value = "select"
